$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 6, pushing the existing rows 6-7 down to 7-8.
$ws.Rows.Item(6).Insert()

# Keep the same custom row height (15.75) the neighbouring rows use.
$ws.Rows.Item(6).RowHeight = 15.75

# Selection moves to reflect the newly inserted area (matches the diff: B13 selected).
$ws.Range("B13").Select()
